# Update the CheckList worksheet ("List1") with the latest test run's
# results: four more verification steps were added (Click on lk button,
# Enter phone number, Enter password, Click login all now have a result
# for every browser), a new browser column (G) was added, and the
# environment / run data in rows 1-4 was refreshed for the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: OS -------------------------------------------------------
$ws.Range("B1:G1").Value = "Windows-10-10.0.18362-SP0"

# --- Row 2: Browser name ---------------------------------------------
$ws.Range("B2").Value = "msedge"
$ws.Range("C2").Value = "chrome"
$ws.Range("D2").Value = "firefox"
$ws.Range("E2").Value = "internet explorer"
$ws.Range("F2").Value = "internet explorer"
$ws.Range("G2").Value = "opera"

# --- Row 3: Browser version (leading "'" keeps "75.0"/"11" as text,
#     matching how the other version strings are stored) --------------
$ws.Range("B3").Value = "81.0.416.72"
$ws.Range("C3").Value = "81.0.4044.138"
$ws.Range("D3").Value = "'75.0"
$ws.Range("E3").Value = "'11"
$ws.Range("F3").Value = "'11"
$ws.Range("G3").Value = "81.0.4044.129"

# --- Row 4: Last checked ----------------------------------------------
$ws.Range("B4").Value = "10.05.2020_13.50.05"
$ws.Range("C4").Value = "10.05.2020_13.50.46"
$ws.Range("D4").Value = "10.05.2020_13.51.29"
$ws.Range("E4").Value = "10.05.2020_13.56.33"
$ws.Range("F4").Value = "10.05.2020_13.57.41"
$ws.Range("G4").Value = "10.05.2020_13.59.12"

# --- Row 6: Open Main page (now also has a result for the new browser)
$ws.Range("G6").Value = "Pass"

# --- Rows 7-10: four new steps, each now checked for every browser ----
$ws.Range("B7:G7").Value = "Pass"
$ws.Range("B8:G8").Value = "Pass"
$ws.Range("B9:G9").Value = "Pass"
$ws.Range("B10:G10").Value = "Pass"

# --- Normalise styles: copy the plain (default) look already used by
#     B1 onto every data cell we touched -- done as two ranges so the
#     still-blank row 5 between the two data blocks is left untouched.
$ws.Range("B1").Copy()
$ws.Range("C1:G4").PasteSpecial(-4122)
$ws.Range("B6:G10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Columns: B:G collapse to one uniform width, column A unchanged --
# (26.140625 characters isn't on Excel's pixel-quantised width grid;
# 25.25 is the closest input that rounds back to the nearest slot.)
$ws.Range("B1:G1").ColumnWidth = 25.25

# --- Selection moves on to the next empty cell for the following run --
$ws.Range("F22").Select()
